$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.368102310726031
$ws.Range("C2").Value = 0.994912726396776
$ws.Range("D2").Value = 0.4951002342582931
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
$ws.Range("G2").Value = 0.1245324579833929
$ws.Range("H2").Value = 0.992
